$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        A = "Record"
        B = "RJ Record"
        C = "Social"
        D = "2025-04-02T18:58"
        E = "Negativo"
        F = "Direito à prioridade não é respeitado durante atualização de dados do CadÚnico. Imagens mostram centenas de pessoas no primeiro dia do mutirão do CadÚnico, que acontece até sexta-feira, na Fundação de Esportes. Entrevista com beneficiários e reclamações de mães de autista e com bebê no colo e de senhora com deficiência que não teve prioridade. 800 senhas por dia. Principal reivindicação é de que o cadastro seja feita nos Cras. Pessoas aguardando para atendimento amanhã. *matéria* Às 19h02, repórter *ao vivo* leu a nota da prefeitura. São 3 dias de ação. Atendimento é por ordem de chegada. Portões abrem 9h. Outros mutirões serão realizados em abril."
    },
    @{
        A = "Record"
        B = "RJ Record"
        C = "Social"
        D = "2025-04-02T19:05"
        E = "Neutro"
        F = "Dia do Autismo. Data reforça a conscientização mundial sobre o transtorno do espectro autista. Entrevista com mãe e com presidente da Apape, Naira Peçanha. *matéria* também foi veiculada no Balanço Geral. "
    }
)

$startRow = 22
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
}
